$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Output_flows")
$ws2 = $wb.Worksheets.Item("Input_flows")

$ws1.Cells.Item(2,3).Value = [double]"7.173791549909988E-30"
$ws1.Cells.Item(2,5).Value = [double]"2.096910438516114E-29"
$ws1.Cells.Item(2,8).Value = [double]"2.092439414126954E-44"
$ws1.Cells.Item(3,3).Value = [double]"1.085305907103524E-30"
$ws1.Cells.Item(3,4).Value = [double]"1.072442089058735E-31"
$ws1.Cells.Item(3,5).Value = [double]"3.172366062988177E-30"
$ws1.Cells.Item(3,8).Value = [double]"3.165601956244065E-45"
$ws1.Cells.Item(4,3).Value = [double]"4.724100676349406E-31"
$ws1.Cells.Item(4,4).Value = [double]"4.668107272896908E-31"
$ws1.Cells.Item(4,5).Value = [double]"1.380861982386761E-30"
$ws1.Cells.Item(4,8).Value = [double]"1.3779177137676E-45"
$ws1.Cells.Item(5,3).Value = [double]"7.840866860464034E-32"
$ws1.Cells.Item(5,4).Value = [double]"7.74793132593293E-31"
$ws1.Cells.Item(5,5).Value = [double]"2.291897590323544E-31"
$ws1.Cells.Item(5,8).Value = [double]"2.287010815077652E-46"
$ws1.Cells.Item(6,3).Value = [double]"2.160068184860726E-31"
$ws1.Cells.Item(6,4).Value = [double]"2.134465519370269E-29"
$ws1.Cells.Item(6,5).Value = [double]"6.313912933249146E-31"
$ws1.Cells.Item(6,8).Value = [double]"6.300450432325373E-46"
$ws1.Cells.Item(7,3).Value = [double]"1.304089766154492E-22"
$ws1.Cells.Item(7,5).Value = [double]"3.811874689118481E-21"
$ws1.Cells.Item(7,8).Value = [double]"3.803747024536139E-36"
$ws1.Cells.Item(8,3).Value = [double]"2.868800764009127E-23"
$ws1.Cells.Item(8,4).Value = [double]"2.835742313334075E-25"
$ws1.Cells.Item(8,5).Value = [double]"8.385549295963613E-22"
$ws1.Cells.Item(8,8).Value = [double]"8.367669659938103E-37"
$ws1.Cells.Item(9,3).Value = [double]"1.906848116239122E-23"
$ws1.Cells.Item(9,4).Value = [double]"1.884247423305785E-24"
$ws1.Cells.Item(9,5).Value = [double]"5.573746730425661E-22"
$ws1.Cells.Item(9,8).Value = [double]"5.561862409038825E-37"
$ws1.Cells.Item(10,3).Value = [double]"9.232744214815683E-24"
$ws1.Cells.Item(10,4).Value = [double]"9.123311162295117E-24"
$ws1.Cells.Item(10,5).Value = [double]"2.698745507937039E-22"
$ws1.Cells.Item(10,8).Value = [double]"2.692991253122674E-37"
$ws1.Cells.Item(11,3).Value = [double]"5.526889724111309E-22"
$ws1.Cells.Item(11,4).Value = [double]"5.461381093505629E-21"
$ws1.Cells.Item(11,5).Value = [double]"1.615518470865224E-20"
$ws1.Cells.Item(11,8).Value = [double]"1.612073868581929E-35"
$ws1.Cells.Item(12,3).Value = [double]"1.775314911890729E-28"
$ws1.Cells.Item(12,5).Value = [double]"2.59463656317397E-27"
$ws1.Cells.Item(12,8).Value = [double]"2.589104289052064E-42"
$ws1.Cells.Item(13,3).Value = [double]"4.210357186149053E-29"
$ws1.Cells.Item(13,4).Value = [double]"3.120339777137902E-29"
$ws1.Cells.Item(13,5).Value = [double]"6.153469801912532E-28"
$ws1.Cells.Item(13,8).Value = [double]"6.140349397217622E-43"
$ws1.Cells.Item(14,3).Value = [double]"2.681835810518359E-29"
$ws1.Cells.Item(14,4).Value = [double]"7.950146351328639E-29"
$ws1.Cells.Item(14,5).Value = [double]"3.919523913078316E-28"
$ws1.Cells.Item(14,8).Value = [double]"3.9111667192338E-43"
$ws1.Cells.Item(15,3).Value = [double]"3.918134701885327E-29"
$ws1.Cells.Item(15,4).Value = [double]"9.872820160999686E-28"
$ws1.Cells.Item(15,5).Value = [double]"5.72638436643636E-28"
$ws1.Cells.Item(15,8).Value = [double]"5.714174591667839E-43"
$ws1.Cells.Item(16,3).Value = [double]"1.215252103903267E-29"
$ws1.Cells.Item(16,4).Value = [double]"3.008124417977028E-27"
$ws1.Cells.Item(16,5).Value = [double]"1.776100409646976E-28"
$ws1.Cells.Item(16,8).Value = [double]"1.772313415170117E-43"
$ws1.Cells.Item(17,3).Value = [double]"5.624630120241189E-21"
$ws1.Cells.Item(17,5).Value = [double]"3.288176281641159E-19"
$ws1.Cells.Item(17,8).Value = [double]"3.28116524479331E-34"
$ws1.Cells.Item(18,3).Value = [double]"1.234518153822326E-21"
$ws1.Cells.Item(18,4).Value = [double]"3.659702478492827E-23"
$ws1.Cells.Item(18,5).Value = [double]"7.217031566299576E-20"
$ws1.Cells.Item(18,8).Value = [double]"7.20164343929254E-35"
$ws1.Cells.Item(19,3).Value = [double]"8.485807342653865E-22"
$ws1.Cells.Item(19,4).Value = [double]"1.006227495843757E-22"
$ws1.Cells.Item(19,5).Value = [double]"4.96082939467926E-20"
$ws1.Cells.Item(19,8).Value = [double]"4.950251933283355E-35"
$ws1.Cells.Item(20,3).Value = [double]"2.433617217588345E-20"
$ws1.Cells.Item(20,4).Value = [double]"2.452867700127614E-20"
$ws1.Cells.Item(20,5).Value = [double]"1.422700202928972E-18"
$ws1.Cells.Item(20,8).Value = [double]"1.419666727016543E-33"
$ws1.Cells.Item(21,3).Value = [double]"4.610601176485249E-20"
$ws1.Cells.Item(21,4).Value = [double]"4.565064955986278E-19"
$ws1.Cells.Item(21,5).Value = [double]"2.695371803750808E-18"
$ws1.Cells.Item(21,8).Value = [double]"2.689624742335561E-33"

$ws2.Cells.Item(2,3).Value = [double]"1.239611184737257E-29"
$ws2.Cells.Item(3,3).Value = [double]"7.748539097786703E-31"
$ws2.Cells.Item(4,3).Value = [double]"3.065863142246586E-32"
$ws2.Cells.Item(5,3).Value = [double]"1.515880054515326E-32"
$ws2.Cells.Item(6,3).Value = [double]"2.219205330551368E-29"
$ws2.Cells.Item(7,3).Value = [double]"1.1261721563963E-22"
$ws2.Cells.Item(8,3).Value = [double]"4.55525614978306E-23"
$ws2.Cells.Item(9,3).Value = [double]"3.036463004524159E-23"
$ws2.Cells.Item(10,3).Value = [double]"1.516155149553344E-23"
$ws2.Cells.Item(11,3).Value = [double]"2.216925477456902E-20"
$ws2.Cells.Item(12,3).Value = [double]"3.157622591123431E-30"
$ws2.Cells.Item(13,3).Value = [double]"2.281521509323375E-31"
$ws2.Cells.Item(14,3).Value = [double]"3.367908605495537E-33"
$ws2.Cells.Item(15,3).Value = [double]"1.448695578863607E-27"
$ws2.Cells.Item(16,3).Value = [double]"3.197886979980758E-27"
$ws2.Cells.Item(17,3).Value = [double]"8.359676497522087E-23"
$ws2.Cells.Item(18,3).Value = [double]"9.409726761523015E-24"
$ws2.Cells.Item(19,3).Value = [double]"1.112470524411943E-24"
$ws2.Cells.Item(20,3).Value = [double]"1.448739727326202E-18"
$ws2.Cells.Item(21,3).Value = [double]"3.197984311114291E-18"
